$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.051558125831509
$ws.Range("D2").Value = 1.053748318504431
$ws.Range("E2").Value = 1.048251097400142
$ws.Range("F2").Value = 1.059505324490593
$ws.Range("I2").Value = 1.049647275507856
$ws.Range("J2").Value = 1.056585365705096
$ws.Range("K2").Value = 1.056493430923028
$ws.Range("L2").Value = 1.051011464325532
$ws.Range("M2").Value = 1.062234645293642
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05294731746347
$ws.Range("D3").Value = 1.054866382000574
$ws.Range("E3").Value = 1.049473553825967
$ws.Range("F3").Value = 1.061085512210101
$ws.Range("I3").Value = 1.050193171645195
$ws.Range("J3").Value = 1.057622546760413
$ws.Range("K3").Value = 1.057423820220613
$ws.Range("L3").Value = 1.052044866240878
$ws.Range("M3").Value = 1.063627139126467
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053844158736299
$ws.Range("D4").Value = 1.055588056272541
$ws.Range("E4").Value = 1.050262215559527
$ws.Range("F4").Value = 1.062106466013311
$ws.Range("I4").Value = 1.05054404292949
$ws.Range("J4").Value = 1.058291167769236
$ws.Range("K4").Value = 1.058023473583255
$ws.Range("L4").Value = 1.052710639686535
$ws.Range("M4").Value = 1.064526154751454
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.054220705838855
$ws.Range("D5").Value = 1.055891025821861
$ws.Range("E5").Value = 1.050593213435465
$ws.Range("F5").Value = 1.062535316866566
$ws.Range("I5").Value = 1.050690987816082
$ws.Range("J5").Value = 1.058571662584957
$ws.Range("K5").Value = 1.058275005659777
$ws.Range("L5").Value = 1.052989841421938
$ws.Range("M5").Value = 1.064903625239757
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.054283901543088
$ws.Range("D6").Value = 1.055941871140526
$ws.Range("E6").Value = 1.050648757059593
$ws.Range("F6").Value = 1.062607302070337
$ws.Range("I6").Value = 1.050715627692171
$ws.Range("J6").Value = 1.05861872427532
$ws.Range("K6").Value = 1.058317206189721
$ws.Range("L6").Value = 1.053036680380432
$ws.Range("M6").Value = 1.064966976537436
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053849192072978
$ws.Range("D7").Value = 1.055592106218937
$ws.Range("E7").Value = 1.05026664054359
$ws.Range("F7").Value = 1.062112197735512
$ws.Range("I7").Value = 1.050546008615753
$ws.Range("J7").Value = 1.058294918079008
$ws.Range("K7").Value = 1.058026836769609
$ws.Range("L7").Value = 1.052714373094771
$ws.Range("M7").Value = 1.064531200389617
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.052028040645141
$ws.Range("D8").Value = 1.054126546829945
$ws.Range("E8").Value = 1.04866472288984
$ws.Range("F8").Value = 1.060039679217146
$ws.Range("I8").Value = 1.049832254490903
$ws.Range("J8").Value = 1.056936408167731
$ws.Range("K8").Value = 1.056808354454157
$ws.Range("L8").Value = 1.051361313085015
$ws.Range("M8").Value = 1.062705668486252
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.0488028187749
$ws.Range("D9").Value = 1.051530099403937
$ws.Range("E9").Value = 1.045823635147013
$ws.Range("F9").Value = 1.056375490799912
$ws.Range("I9").Value = 1.048556291757873
$ws.Range("J9").Value = 1.054523083058827
$ws.Range("K9").Value = 1.054642834580095
$ws.Range("L9").Value = 1.048954496490033
$ws.Range("M9").Value = 1.059473002969653
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.046641335068542
$ws.Range("D10").Value = 1.049789401259322
$ws.Range("E10").Value = 1.043916832333014
$ws.Range("F10").Value = 1.053923966951429
$ws.Range("I10").Value = 1.047693160610111
$ws.Range("J10").Value = 1.052900729974029
$ws.Range("K10").Value = 1.053186444732701
$ws.Range("L10").Value = 1.047334381436358
$ws.Range("M10").Value = 1.057306724343322
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045702588659877
$ws.Range("D11").Value = 1.049033270768514
$ws.Range("E11").Value = 1.043088044765381
$ws.Range("F11").Value = 1.052860232595178
$ws.Range("I11").Value = 1.047316401768258
$ws.Range("J11").Value = 1.052194951347476
$ws.Range("K11").Value = 1.052552720887712
$ws.Range("L11").Value = 1.046629069129724
$ws.Range("M11").Value = 1.056365937556086
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.045353464631603
$ws.Range("D12").Value = 1.048752043138157
$ws.Range("E12").Value = 1.042779717460563
$ws.Range("F12").Value = 1.05246477153917
$ws.Range("I12").Value = 1.047175998949095
$ws.Range("J12").Value = 1.051932292012034
$ws.Range("K12").Value = 1.052316855613463
$ws.Range("L12").Value = 1.046366507086502
$ws.Range("M12").Value = 1.056016060701022
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045428372617767
$ws.Range("D13").Value = 1.048812384144973
$ws.Range("E13").Value = 1.042845876514343
$ws.Range("F13").Value = 1.052549614965606
$ws.Range("I13").Value = 1.047206136619745
$ws.Range("J13").Value = 1.051988656157945
$ws.Range("K13").Value = 1.052367470983133
$ws.Range("L13").Value = 1.04642285380463
$ws.Range("M13").Value = 1.056091129947336
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045673738824016
$ws.Range("D14").Value = 1.049010031934424
$ws.Range("E14").Value = 1.043062568149238
$ws.Range("F14").Value = 1.052827550716798
$ws.Range("I14").Value = 1.047304805392709
$ws.Range("J14").Value = 1.052173250124122
$ws.Range("K14").Value = 1.05253323385845
$ws.Range("L14").Value = 1.046607377499547
$ws.Range("M14").Value = 1.056337025389986
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045824859530915
$ws.Range("D15").Value = 1.049131760405309
$ws.Range("E15").Value = 1.04319601532777
$ws.Range("F15").Value = 1.052998750358127
$ws.Range("I15").Value = 1.047365537681034
$ws.Range("J15").Value = 1.05228691784696
$ws.Range("K15").Value = 1.052635303088834
$ws.Range("L15").Value = 1.046720991854295
$ws.Range("M15").Value = 1.056488472837697
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.046703576553772
$ws.Range("D16").Value = 1.049839532056078
$ws.Range("E16").Value = 1.043971769574212
$ws.Range("F16").Value = 1.053994516035382
$ws.Range("I16").Value = 1.047718100899351
$ws.Range("J16").Value = 1.052947500234684
$ws.Range("K16").Value = 1.053228437070133
$ws.Range("L16").Value = 1.047381110123394
$ws.Range("M16").Value = 1.057369101969378
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.047254013476258
$ws.Range("D17").Value = 1.050282851910469
$ws.Range("E17").Value = 1.044457535885567
$ws.Range("F17").Value = 1.054618534580412
$ws.Range("I17").Value = 1.04793844350444
$ws.Range("J17").Value = 1.053360980209706
$ws.Range("K17").Value = 1.053599660387693
$ws.Range("L17").Value = 1.0477941640935
$ws.Range("M17").Value = 1.057920747676792
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.047574803584363
$ws.Range("D18").Value = 1.050541202111872
$ws.Range("E18").Value = 1.044740573799724
$ws.Range("F18").Value = 1.054982301712596
$ws.Range("I18").Value = 1.048066674792103
$ws.Range("J18").Value = 1.053601839255904
$ws.Range("K18").Value = 1.053815890265078
$ws.Range("L18").Value = 1.048034726175004
$ws.Range("M18").Value = 1.058242246330083
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04768413914225
$ws.Range("D19").Value = 1.050629253922889
$ws.Range("E19").Value = 1.044837031556053
$ws.Range("F19").Value = 1.055106301220021
$ws.Range("I19").Value = 1.04811034915324
$ws.Range("J19").Value = 1.053683912455191
$ws.Range("K19").Value = 1.053889568651643
$ws.Range("L19").Value = 1.048116689905063
$ws.Range("M19").Value = 1.05835182408392
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047194984846377
$ws.Range("D20").Value = 1.050235311806029
$ws.Range("E20").Value = 1.044405448990245
$ws.Range("F20").Value = 1.054551605337648
$ws.Range("I20").Value = 1.047914832942193
$ws.Range("J20").Value = 1.053316650547438
$ws.Range("K20").Value = 1.053559862552458
$ws.Range("L20").Value = 1.047749885155167
$ws.Range("M20").Value = 1.05786158899798
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045601496590286
$ws.Range("D21").Value = 1.048951839777125
$ws.Range("E21").Value = 1.042998771147156
$ws.Range("F21").Value = 1.0527457151157
$ws.Range("I21").Value = 1.047275762572354
$ws.Range("J21").Value = 1.052118905713412
$ws.Range("K21").Value = 1.052484433907159
$ws.Range("L21").Value = 1.046553055880829
$ws.Range("M21").Value = 1.056264627172255
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04459710407348
$ws.Range("D22").Value = 1.048142742522504
$ws.Range("E22").Value = 1.042111563055202
$ws.Range("F22").Value = 1.051608291689326
$ws.Range("I22").Value = 1.046871302666645
$ws.Range("J22").Value = 1.051362929342584
$ws.Range("K22").Value = 1.051805534915628
$ws.Range("L22").Value = 1.045797215718906
$ws.Range("M22").Value = 1.055258079800747
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04512979208108
$ws.Range("D23").Value = 1.048571864304035
$ws.Range("E23").Value = 1.042582154657735
$ws.Range("F23").Value = 1.052211453548913
$ws.Range("I23").Value = 1.047085967307329
$ws.Range("J23").Value = 1.05176396478807
$ws.Range("K23").Value = 1.052165693539449
$ws.Range("L23").Value = 1.046198220740169
$ws.Range("M23").Value = 1.055791907487089
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047221658179144
$ws.Range("D24").Value = 1.050256793846135
$ws.Range("E24").Value = 1.044428985747514
$ws.Range("F24").Value = 1.054581848437736
$ws.Range("I24").Value = 1.047925502437678
$ws.Range("J24").Value = 1.053336682192775
$ws.Range("K24").Value = 1.053577846403301
$ws.Range("L24").Value = 1.047769894030337
$ws.Range("M24").Value = 1.05788832108348
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.049638577237199
$ws.Range("D25").Value = 1.052203032421536
$ws.Range("E25").Value = 1.04656033839806
$ws.Range("F25").Value = 1.05732426862
$ws.Range("I25").Value = 1.048888343203396
$ws.Range("J25").Value = 1.055149330934432
$ws.Range("K25").Value = 1.055204888483482
$ws.Range("L25").Value = 1.049579430592632
$ws.Range("M25").Value = 1.060310654973627
